$d = $word.ActiveDocument

$pairs = @(
    @{old="196÷9="; new="321÷2="},
    @{old="947÷7="; new="536÷2="},
    @{old="226÷2="; new="517÷6="},
    @{old="370÷7="; new="275÷9="},
    @{old="916÷4="; new="870÷3="},
    @{old="914÷6="; new="637÷2="},
    @{old="229÷8="; new="963÷2="},
    @{old="681÷8="; new="119÷3="},
    @{old="382÷5="; new="372÷8="},
    @{old="231÷4="; new="699÷8="},
    @{old="576÷6="; new="654÷6="},
    @{old="334÷7="; new="242÷8="},
    @{old="559÷3="; new="236÷6="},
    @{old="271÷7="; new="397÷5="},
    @{old="249÷9="; new="495÷2="},
    @{old="846÷5="; new="963÷2="},
    @{old="570÷7="; new="879÷7="},
    @{old="972÷3="; new="555÷8="},
    @{old="341÷8="; new="638÷4="},
    @{old="172÷4="; new="450÷6="},
    @{old="606÷9="; new="109÷8="},
    @{old="634÷7="; new="874÷6="},
    @{old="993÷3="; new="114÷5="},
    @{old="936÷4="; new="393÷9="},
    @{old="863÷5="; new="128÷9="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
